$d = $word.ActiveDocument

# 1) Collapse paragraph 1's three runs (incl. the "alt" spell-check exception
#    wrapped in proofErr markers) into a single plain run with the same text.
$d.Content.Find.Execute(
    "Elijo hacer el alt de las imágenes en la sección de galería con el párrafo para minimizar las cosas a editar cuando se cambien los artistas expuestos",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Elijo hacer el alt de las imágenes en la sección de galería con el párrafo para minimizar las cosas a editar cuando se cambien los artistas expuestos",
    2
)

# 2) Collapse paragraph 2's three runs (incl. the "ese mimos" grammar-check
#    exception wrapped in proofErr markers) into a single plain run.
$d.Content.Find.Execute(
    "Por ese mimos motivo hago una sola “ficha” que va rellenándose en función de la imagen pulsada, y edito desde java el párrafo en función a un párrafo no visible junto al nombre del artistas.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Por ese mimos motivo hago una sola “ficha” que va rellenándose en función de la imagen pulsada, y edito desde java el párrafo en función a un párrafo no visible junto al nombre del artistas.",
    2
)

# 3) Append a new, underlined paragraph at the very end of the document.
$endRange = $d.Content
$endRange.Collapse(0)

$newParaXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:t>Decido que no quiero que las barras se carguen al hacer scroll, si no que comiencen a cargarse a partir de entrar en visión por el usuario.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$endRange.InsertXML($newParaXml)
